$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, pushing old rows 91-93 down to 92-94
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new data record
$ws.Range("A91").Value = 6
$ws.Range("B91").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C91").Value = "Metropolitana"
$ws.Range("D91").Value = 44448
$ws.Range("E91").Value = 13
$ws.Range("F91").Value = 100112001
$ws.Range("G91").Value = "Berenjena"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 280
$ws.Range("K91").Value = 7000
$ws.Range("L91").Value = 8000
$ws.Range("M91").Value = 7571
$ws.Range("N91").Value = "$/caja 50 unidades"
$ws.Range("O91").Value = "Región de Arica y Parinacota"
$ws.Range("P91").Value = 151
$ws.Range("Q91").Value = 50
$ws.Range("R91").Value = "Hortaliza"

# Ensure the date cell keeps the same date number format as the column
$ws.Range("D91").NumberFormat = $ws.Range("D92").NumberFormat
